# Configurable Proto buff example added
#
# Updates the "CreateOrder-OrderCreated-Event" sheet so the REST call and
# the KAFKA verification step use the new "proto" naming / endpoint instead
# of the old JSON-based example, and clears the now-unused JSON request
# body / response-mapping / path-mapping cells on row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CREATE_ORDER / REST step)
$ws.Range("C2").Value = "orderservice-proto"
$ws.Range("F2").Value = "http://localhost:8800/demo/10/Elan"
$ws.Range("J2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("N2").Value = 200
$ws.Range("O2").ClearContents()

# Row 3 (VERIFY_ORDER_CREATED_EVENT / KAFKA step)
$ws.Range("C3").Value = "order-proto"
$ws.Range("K3").Value = "ProtoBuffMessageType"

# Update the active selection saved with the sheet view (also drops the
# previous topLeftCell scroll position, matching the sheet being scrolled
# back to the default top-left).
[void]$ws.Range("C2").Select()
